$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of contact data (write order matches the shared-string
# insertion order: Last Name, Personal Email, First Name).
$ws.Range("B3").Value = "poonacha"
$ws.Range("C3").Value = "terry.poonacha@gmail.com"
$ws.Range("A3").Value = "Terry"

# Turn the new email address into a mailto hyperlink, then restore the
# built-in "Hyperlink" cell style so the style index isn't duplicated.
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:terry.poonacha@gmail.com")
$ws.Cells.Item(3, 3).Style = "Hyperlink"

# Match the final selection left behind in the saved workbook.
[void]$ws.Range("A6").Select()
